# Automatische test-sync: 2025-08-06 19:46:50
# Append a new log row to the "Logs" sheet and bump the matching
# "Dashboard" category counter.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 8 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(8, 1).Value = "Bestel je 100 M5-bouten zodra je kan?"
$logs.Cells.Item(8, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(8, 3).Value = "Testmail #1: Bestel je 100 M5-bouten zodra je kan?"
$logs.Cells.Item(8, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(8, 5).Value = "Beste klant,`nDank je wel voor je bericht. Als je 100 M5-bouten wilt bestellen, kun je dit doen door naar onze website te gaan en het gewenste aantal toe te voegen aan je winkelwagen. Als je hulp nodig hebt bij het plaatsen van de bestelling, laat het ons dan weten en we helpen je graag verder.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(8, 6).Value = "2025-08-06 19:46:25"
$logs.Cells.Item(8, 7).Value = "Ja"
$logs.Cells.Item(8, 8).Value = "Nee"
$logs.Cells.Item(8, 9).Value = "Ja"
$logs.Cells.Item(8, 10).Value = "Nee"

# Undo the engine's automatic row-height bump from the multi-line text in
# column E so row 8 keeps the sheet's default (unmodified) row height.
$logs.Rows.Item(8).AutoFit()

# --- Logs sheet: extend the conditional-formatting ranges to row 8 ------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`7")
    $newRange = $logs.Range("$col`2:$col`8")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: bump the "Inkoop / Bestellingen" tally ------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 2
